$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete row 555 entirely ("「彼は応えるだろう」..." post), shifting all rows below up by one.
$ws.Rows.Item(555).Delete()
